$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage so literal strings (dates, numbers-as-text) are not
# auto-converted by Excel's type inference.
$textRange = $ws.Range("A1:K2")
$textRange.NumberFormat = "@"

# Row 1
$ws.Range("A1").Value = "Ernsberger"
$ws.Range("B1").Value = "Donnie"
$ws.Range("C1").Value = "TE"
$ws.Range("D1").Value = "2018-12-23"
$ws.Range("E1").Value = "15"
$ws.Range("F1").Value = "22.071"
$ws.Range("G1").Value = "TAM"
$ws.Range("H1").Value = "@"
$ws.Range("I1").Value = "DAL"
$ws.Range("J1").Value = "L 20-27"
$ws.Range("K1").Value = ""
$ws.Range("L1").Value = 0

# Row 2
$ws.Range("A2").Value = "Ernsberger"
$ws.Range("B2").Value = "Donnie"
$ws.Range("C2").Value = "TE"
$ws.Range("D2").Value = "2018-12-30"
$ws.Range("E2").Value = "16"
$ws.Range("F2").Value = "22.078"
$ws.Range("G2").Value = "TAM"
$ws.Range("H2").Value = ""
$ws.Range("I2").Value = "ATL"
$ws.Range("J2").Value = "L 32-34"
$ws.Range("K2").Value = ""
$ws.Range("L2").Value = 0
